$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.054833169165874551,
    -0.0099999995808097708,
    -0.0089999995854928017,
    0.28399259543883559,
    -0.0059999995979014287,
    0.033877724115711771,
    -0.019999999505142085,
    -0.042169433656490618,
    -0.0059999995791919503,
    -0.0059999995801334194,
    0.02393596484794358,
    -0.0059999995800175121,
    -0.0059999995798811767,
    0.032802606050978511,
    -0.0059999995785773308,
    -0.0059999995771491399,
    -0.0059999995752901825,
    -0.0089999995585055004,
    -0.008999999585962204,
    -0.0089999995823148993,
    -0.0089999995817304779,
    -0.0089999995813174749,
    -0.0089999995817802159,
    -0.041999999395549992,
    -0.041999999392253073,
    -0.066155039994889364,
    -0.0059999995770922965,
    -0.0059999995679556051,
    -0.011999999528832461,
    -0.019999999482133379,
    -0.014999999504789585,
    -0.020999999471387198,
    -0.0059999995528823291
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
